$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the shared-string text values (the author shortened / fixed two labels)
$ws.Range("B8").Value = "Studien England"
$ws.Range("B69").Value = "Badenische Wahlreform "

# Update the sheet's on-screen selection to match the saved view
$ws.Range("B69").Select()
